$d = $word.ActiveDocument

$replacements = @(
    @("53÷6=8, 5", "49÷2=24, 1"),
    @("85÷5=17, 0", "58÷4=14, 2"),
    @("49÷5=9, 4", "62÷3=20, 2"),
    @("35÷8=4, 3", "77÷3=25, 2"),
    @("84÷6=14, 0", "28÷6=4, 4"),
    @("21÷2=10, 1", "14÷3=4, 2"),
    @("69÷5=13, 4", "63÷6=10, 3"),
    @("52÷6=8, 4", "23÷2=11, 1"),
    @("80÷7=11, 3", "54÷4=13, 2"),
    @("32÷9=3, 5", "64÷9=7, 1"),
    @("73÷4=18, 1", "15÷5=3, 0"),
    @("89÷7=12, 5", "15÷4=3, 3"),
    @("54÷3=18, 0", "28÷4=7, 0"),
    @("38÷5=7, 3", "28÷2=14, 0"),
    @("99÷9=11, 0", "19÷5=3, 4"),
    @("85÷9=9, 4", "30÷2=15, 0"),
    @("62÷8=7, 6", "83÷9=9, 2"),
    @("93÷5=18, 3", "51÷5=10, 1"),
    @("75÷2=37, 1", "76÷2=38, 0"),
    @("48÷5=9, 3", "29÷4=7, 1"),
    @("14÷5=2, 4", "50÷2=25, 0"),
    @("61÷3=20, 1", "69÷3=23, 0"),
    @("80÷4=20, 0", "38÷7=5, 3"),
    @("11÷8=1, 3", "39÷4=9, 3"),
    @("58÷7=8, 2", "40÷2=20, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
